# edit.ps1 -- applies the Kick-Off.docx revision described in the diff
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: center align "Kick-Off do Projeto"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.Alignment = 1   # wdAlignParagraphCenter

# ---------------------------------------------------------------------------
# 2) "Projeto: ..." paragraph -> bold
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 3) "Cliente: ..." paragraph -> bold
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 4) Define the "Sumrio1" (TOC 1) paragraph style used later on, so the
#    upcoming InsertXML call can reference w:pStyle val="Sumrio1".
# ---------------------------------------------------------------------------
$sumrio1 = $d.Styles.Add("Sumrio1", 1)   # wdStyleTypeParagraph
$sumrio1.NameLocal = "toc 1"
$sumrio1.BaseStyle = $d.Styles.Item("Normal")
$sumrio1.NextParagraphStyle = $d.Styles.Item("Normal")
$sumrio1.Priority = 39
$sumrio1.UnhideWhenUsed = $true

# ---------------------------------------------------------------------------
# 5) Replace the whole "Produtos ... Proximos Passos" section (old
#    paragraphs 12-40) with the new termo-de-abertura / TOC block.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(12)
$pEnd = $d.Paragraphs.Item(40)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$newBlockXml = @'
<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p>
      <w:r>
        <w:t>Com todas as informações colhidas, o gerente de projetos elabora o termo de abertura indicando:</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Gerente de projeto</w:t>
      </w:r>
      <w:r>
        <w:t>: John</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Data de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>inicio</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> do projeto e suas </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">dependências: 25,26 do mês 03, na reunião de iniciação é necessário estar presentes todas as partes interessadas, como </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CEO’s</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, CFO, CIO, CTO, COO.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Requisitos que satisfazem as necessidades do cliente</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t>Elaborar um software para auxiliar na gestão de viagens e logística.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Justificativa do projeto</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: resolver o problema do cliente, que seria </w:t>
      </w:r>
      <w:r>
        <w:t>o gerenciamento de viagens, gestão de logística.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>WBS</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Tomar conhecimento sobre o problema do cliente, propondo uma solução para o mesmo, estipular datas de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>inicio</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, dependências e um orçamento sobre o que será proposto.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio1"/>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="539"/>
          <w:tab w:val="left" w:pos="702"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9702"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:spacing w:before="100" w:after="0" w:line="277" w:lineRule="exact"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:hyperlink w:anchor="_TOC_250047" w:history="1">
        <w:r>
          <w:rPr>
            <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
          <w:t>FASES DO PROJETO</w:t>
        </w:r>
      </w:hyperlink>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio2"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1321"/>
          <w:tab w:val="left" w:pos="1322"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
        <w:spacing w:before="0" w:line="277" w:lineRule="exact"/>
      </w:pPr>
      <w:hyperlink w:anchor="_TOC_250046" w:history="1">
        <w:r>
          <w:t>Definição</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:spacing w:val="-1"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:t>de</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:spacing w:val="-1"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r>
          <w:t>Requisitos</w:t>
        </w:r>
      </w:hyperlink>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio3"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1541"/>
          <w:tab w:val="left" w:pos="1542"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
        <w:spacing w:before="103"/>
      </w:pPr>
      <w:hyperlink w:anchor="_TOC_250043" w:history="1">
        <w:r>
          <w:t>Artefatos da Concepção</w:t>
        </w:r>
      </w:hyperlink>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:r>
        <w:t>Termo de Abertura do Projeto</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Reunião de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>K</w:t>
      </w:r>
      <w:r>
        <w:t>ick</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>-Off</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:r>
        <w:t>Briefing</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:r>
        <w:t>Especificação Técnica do Projeto</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:r>
        <w:t>Business Case</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="Sumrio4"/>
        <w:numPr>
          <w:ilvl w:val="3"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="1981"/>
          <w:tab w:val="left" w:pos="1982"/>
          <w:tab w:val="right" w:leader="dot" w:pos="9285"/>
        </w:tabs>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:t>ermo de Aceite</w:t>
      </w:r>
    </w:p><w:p/><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Orçamento resumido</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">R$ </w:t>
      </w:r>
      <w:r>
        <w:t>30 mil</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>( trinta</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> mil reais ).</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:pStyle w:val="PargrafodaLista"/>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($newBlockXml)

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
